$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new column S (year 2022) by copying formatting from column R ---

# Row 2 (empty, formatted cell only)
$null = $ws.Range("R2").Copy()
$null = $ws.Range("S2").PasteSpecial(-4122)

# Row 3 header (year 2022)
$null = $ws.Range("R3").Copy()
$null = $ws.Range("S3").PasteSpecial(-4122)
$ws.Range("S3").Value = 2022

# Row 4 data
$null = $ws.Range("R4").Copy()
$null = $ws.Range("S4").PasteSpecial(-4122)
$ws.Range("S4").Value = 13.6

# Row 5 data
$null = $ws.Range("R5").Copy()
$null = $ws.Range("S5").PasteSpecial(-4122)
$ws.Range("S5").Value = 20

# --- Update existing values in row 4 ---
$ws.Range("P4").Value = 13.7
$ws.Range("Q4").Value = 13.1
$ws.Range("R4").Value = 11.8

# --- Update existing values in row 5 ---
$ws.Range("P5").Value = 13.6
$ws.Range("Q5").Value = 12.5
$ws.Range("R5").Value = 13.5

# --- Update selection to S2 ---
[void]$ws.Range("S2").Select()
